$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1378.98894257555
$ws.Range("B3").Value = 1824.012946008359
$ws.Range("B4").Value = 2261.361683477984
$ws.Range("B5").Value = 2734.894384922172
$ws.Range("B6").Value = 3138.72843865085
$ws.Range("B7").Value = 3221.372219969964
$ws.Range("B8").Value = 3435.668064787017
$ws.Range("B9").Value = 3563.320758137213
$ws.Range("B10").Value = 3784.01931899895
$ws.Range("B11").Value = 3918.280102258143
$ws.Range("B12").Value = 4044.368272677009
$ws.Range("B13").Value = 4181.71705606111
$ws.Range("B14").Value = 4258.991484663136
$ws.Range("B15").Value = 4366.326199333155
$ws.Range("B16").Value = 4418.783357305564
$ws.Range("B17").Value = 4522.954156992219
$ws.Range("B18").Value = 4549.649989672791
$ws.Range("B19").Value = 4576.360328042001
$ws.Range("B20").Value = 4627.261578260298
$ws.Range("B21").Value = 4634.10746737402
$ws.Range("B22").Value = 4668.328395554517
$ws.Range("B23").Value = 4651.391857517373
$ws.Range("B24").Value = 4652.667294079483
$ws.Range("B25").Value = 4640.685488817128
$ws.Range("B26").Value = 4616.307397768557
$ws.Range("B27").Value = 4568.980455365257
$ws.Range("B28").Value = 4511.688148036834
$ws.Range("B29").Value = 4470.798258772613
$ws.Range("B30").Value = 4422.922476553295
$ws.Range("B31").Value = 4352.386711579516
$ws.Range("B32").Value = 4306.678444237267
$ws.Range("B33").Value = 4247.147152509612
$ws.Range("B34").Value = 4148.517981610437
$ws.Range("B35").Value = 4075.040902889667
$ws.Range("B36").Value = 3972.652728562359
$ws.Range("B37").Value = 3916.117878662214
$ws.Range("B38").Value = 3800.484964768681
$ws.Range("B39").Value = 3691.025746438188
$ws.Range("B40").Value = 3585.448713558844
$ws.Range("B41").Value = 3501.732628776179
$ws.Range("B42").Value = 3326.428267445979
$ws.Range("B43").Value = 3201.689281562613
$ws.Range("B44").Value = 3063.171019903505
$ws.Range("B45").Value = 2965.046458180914
$ws.Range("B46").Value = 2933.17071835558
$ws.Range("B47").Value = 2829.794468065788
$ws.Range("B48").Value = 2646.460991591533
$ws.Range("B49").Value = 2494.669511529986
$ws.Range("B50").Value = 2340.907175172379
$ws.Range("B51").Value = 2241.671840041132
$ws.Range("B52").Value = 2063.436213330111
$ws.Range("B53").Value = 1826.080387561454
$ws.Range("B54").Value = 1675.26645715361
$ws.Range("B55").Value = 1539.850758386492
$ws.Range("B56").Value = 1454.219119239909
$ws.Range("B57").Value = 1245.247695009602
$ws.Range("B58").Value = 1148.673334318132
$ws.Range("B59").Value = 1066.238904734081
$ws.Range("B60").Value = 1044.486566158677
$ws.Range("B61").Value = 1020.815138994861
$ws.Range("B62").Value = 992.3528348118001
